# Update "想去人数" (want-to-go count) values for a couple of events
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1098
$ws1.Range("F4").Value = 1747
$ws1.Range("F6").Value = 77

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1098
$ws4.Range("F4").Value = 1747
$ws4.Range("F7").Value = 77
